# Updated cryptos list on Mon May  8 08:15:18 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell([string]$addr, [string]$val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "28.020.07"
$ws.Range("E2").Value = "  -3.37%  "

# Row 3
Set-TextCell "D3" "1.857.68"
$ws.Range("E3").Value = "  -2.91%  "

# Row 4
Set-TextCell "D4" "1.006"
$ws.Range("E4").Value = "  +0.44%  "

# Row 5
Set-TextCell "D5" "317.57"
$ws.Range("E5").Value = "  -2.29%  "

# Row 6
$ws.Range("E6").Value = "  +0.44%  "

# Row 7
Set-TextCell "D7" "0.4366"
$ws.Range("E7").Value = "  -4.90%  "

# Row 8
Set-TextCell "D8" "0.3668"
$ws.Range("E8").Value = "  -4.10%  "

# Row 9
Set-TextCell "D9" "0.07472"
$ws.Range("E9").Value = "  -3.20%  "

# Row 10
Set-TextCell "D10" "0.9320"
$ws.Range("E10").Value = "  -5.07%  "

# Row 11
Set-TextCell "D11" "21.26"
$ws.Range("E11").Value = "  -3.89%  "

# Row 12
Set-TextCell "D12" "1.877.31"
$ws.Range("E12").Value = "  -1.00%  "

# Row 13
Set-TextCell "D13" "6.676"
$ws.Range("E13").Value = "  -3.94%  "

# Row 14
Set-TextCell "D14" "5.414"

# Row 15
Set-TextCell "D15" "0.06892"
$ws.Range("E15").Value = "  -1.94%  "

# Row 16
Set-TextCell "D16" "1.007"
$ws.Range("E16").Value = "  +0.45%  "

# Row 17
Set-TextCell "D17" "81.41"
$ws.Range("E17").Value = "  -3.09%  "

# Row 18
Set-TextCell "D18" "0.000008976"
$ws.Range("E18").Value = "  -5.24%  "

# Row 19
Set-TextCell "D19" "1.004"
$ws.Range("E19").Value = "  +0.32%  "

# Row 20
Set-TextCell "D20" "15.77"
$ws.Range("E20").Value = "  -5.82%  "

# Row 21
Set-TextCell "D21" "28.034.49"
$ws.Range("E21").Value = "  -3.23%  "

# Row 22
Set-TextCell "D22" "5.099"
$ws.Range("E22").Value = "  -4.19%  "

# Row 23
$ws.Range("E23").Value = "  -1.33%  "

# Row 24
Set-TextCell "D24" "2.116.53"
$ws.Range("E24").Value = "  -0.44%  "

# Row 25
Set-TextCell "D25" "1.999"
$ws.Range("E25").Value = "  -4.28%  "

# Row 26
Set-TextCell "D26" "154.24"
$ws.Range("E26").Value = "  -2.82%  "

# Row 27
$ws.Range("E27").Value = "  -3.46%  "

# Row 28
Set-TextCell "D28" "5.307"
$ws.Range("E28").Value = "  -6.56%  "

# Row 29
Set-TextCell "D29" "113.17"
$ws.Range("E29").Value = "  -3.67%  "

# Row 30
Set-TextCell "D30" "1.723"
$ws.Range("E30").Value = "  -7.34%  "

# Row 31
Set-TextCell "D31" "0.08981"

# Row 32
Set-TextCell "D32" "0.7925"
$ws.Range("E32").Value = "  -8.62%  "

# Row 33
Set-TextCell "D33" "4.817"
$ws.Range("E33").Value = "  -5.26%  "

# Row 34
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell "D34" "3.035"
$ws.Range("E34").Value = "  +0.34%  "

# Row 35
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D35" "1.166"
$ws.Range("E35").Value = "  -7.10%  "

# Row 36
Set-TextCell "D36" "1.004"
$ws.Range("E36").Value = "  +0.36%  "

# Row 37
Set-TextCell "D37" "1.120"
$ws.Range("E37").Value = "  -3.01%  "

# Row 38
Set-TextCell "D38" "0.05383"
$ws.Range("E38").Value = "  -6.18%  "

# Row 40
Set-TextCell "D40" "2.922"
$ws.Range("E40").Value = "  +2.01%  "

# Row 41
Set-TextCell "D41" "0.5228"
$ws.Range("E41").Value = "  -5.12%  "

# Row 42
Set-TextCell "D42" "6.965"

# Row 43
Set-TextCell "D43" "0.1670"
$ws.Range("E43").Value = "  -4.98%  "

# Row 44
Set-TextCell "D44" "8.716"
$ws.Range("E44").Value = "  -6.61%  "

# Row 45
Set-TextCell "D45" "0.06713"
$ws.Range("E45").Value = "  -2.25%  "

# Row 46
Set-TextCell "D46" "0.4842"
$ws.Range("E46").Value = "  -6.59%  "

# Row 47
Set-TextCell "D47" "10.58"
$ws.Range("E47").Value = "  -6.05%  "

# Row 48
Set-TextCell "D48" "106.82"
$ws.Range("E48").Value = "  -3.40%  "

# Row 49
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell "D49" "1.004"
$ws.Range("E49").Value = "  +0.34%  "

# Row 50
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D50" "1.903"
$ws.Range("E50").Value = "  -7.54%  "

# Row 51
Set-TextCell "D51" "1.663"
$ws.Range("E51").Value = "  -6.73%  "
